# Update countries & provincias Spain
# This script refreshes the COVID data table on sheet "Pais": it updates
# the "last updated" timestamp in A1, updates the statistics (Casos
# totales, Nuevos casos, Casos activos, Recuperados, Casos criticos,
# Muertes hoy, Muertes) for every country whose numbers changed, and
# re-orders a handful of rows whose countries swapped rank (the country
# name in column A and its row data move together).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp ------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 18 de Septiembre de 2020 a las 17:55"

# --- Helper to write a full data row (A..H) ------------------------------
function Set-Row($Row, $Country, $Total, $Nuevos, $Activos, $Recuperados, $Criticos, $MuertesHoy, $Muertes) {
    $ws.Cells.Item($Row, 1).Value = $Country
    $ws.Cells.Item($Row, 2).Value = $Total
    $ws.Cells.Item($Row, 3).Value = $Nuevos
    $ws.Cells.Item($Row, 4).Value = $Activos
    $ws.Cells.Item($Row, 5).Value = $Recuperados
    $ws.Cells.Item($Row, 6).Value = $Criticos
    $ws.Cells.Item($Row, 7).Value = $MuertesHoy
    $ws.Cells.Item($Row, 8).Value = $Muertes
}

# Row 4: Estados Unidos
Set-Row 4 "Estados Unidos" 6885289 10693 4157714 2525196 0 166 202379

# Row 5: India
Set-Row 5 "India" 5247627 34941 4143247 1019663 0 313 84717

# Row 12: España
Set-Row 12 "España" 654637 0 0 0 0 0 30405

# Row 17: Reino Unido
Set-Row 17 "Reino Unido" 385936 4322 0 0 0 27 41732

# Row 20: Irak
Set-Row 20 "Irak" 311690 4305 245305 57977 0 76 8408

# Row 23: Italia
Set-Row 23 "Italia" 294932 1907 216807 42457 0 10 35668

# Row 25: Alemania
Set-Row 25 "Alemania" 270653 1611 241300 19890 0 6 9463

# Row 29: Canada
Set-Row 29 "Canada" 141268 401 123304 8764 0 0 9200

# Row 33: Rumania
Set-Row 33 "Rumania" 110217 1527 88235 17622 0 48 4360

# Row 46: Emiratos Arabes Unidos
Set-Row 46 "Emiratos Arabes Unidos" 83433 865 72790 10240 0 1 403

# Row 57: Singapur
Set-Row 57 "Singapur" 57543 11 57071 445 0 0 27

# Row 60: Suiza
Set-Row 60 "Suiza" 49283 488 39900 7338 0 3 2045

# Rows 64-65: Moldavia overtakes Kirguistan
Set-Row 64 "Moldavia" 45648 665 33734 10728 0 16 1186
Set-Row 65 "Kirguistan" 45244 91 41415 2766 0 0 1063

# Rows 88-90: Grecia overtakes Senegal and Croacia
Set-Row 88 "Grecia" 14738 338 3804 10607 0 2 327
Set-Row 89 "Senegal" 14645 27 11051 3293 0 1 301
Set-Row 90 "Croacia" 14513 234 12169 2100 0 6 244

# Row 94: Albania
Set-Row 94 "Albania" 12073 125 6831 4889 0 6 353

# Row 99: Guayana Francesa
Set-Row 99 "Guayana Francesa" 9659 36 9298 296 0 0 65

# Rows 127-128: Jordania overtakes Eslovenia
Set-Row 127 "Jordania" 4344 213 2511 1804 0 3 29
Set-Row 128 "Eslovenia" 4195 137 2939 1116 0 4 140

# Row 159: Republica de Chipre
Set-Row 159 "Republica de Chipre" 1565 7 1282 261 0 0 22

# Rows 214-215: Islas Malvinas overtakes Montserrat
Set-Row 214 "Islas Malvinas" 13 0 13 0 0 0 0
Set-Row 215 "Montserrat" 13 0 12 0 0 0 1
